$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A14").Value = "24/10/2025"
$ws.Range("B14").Value = "Termalica B-B."
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = "Zaglebie"
$ws.Range("F14").Value = "D"
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 1
$ws.Range("K14").Value = 1.43
$ws.Range("L14").Value = 1.08
$ws.Range("M14").Value = 11
$ws.Range("N14").Value = 13
$ws.Range("O14").Value = 3
$ws.Range("P14").Value = 6
